$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'Edn1', 'Ednrb', 'ECs', 3, 1, 20.66988266666667, 62.009648, 0.9592026714402365, 0.9592026714402365, 3, 1, 20.92735966666666, 62.782079, 0.717329969634113, 0.717329969634113, 432.5660688331324, 3893.094619498192, 0.688064823177185, 0.688064823177185),
    @('ECs', 'Edn1', 'Ednrb', 'FAPs', 3, 1, 20.66988266666667, 62.009648, 0.9592026714402365, 0.9592026714402365, 1, 0.3333333333333333, 0.283297, 0.849891, 0.009710610016949358, 0.009710610016949358, 5.855715749818666, 52.701441748368, 0.009314443069572145, 0.009314443069572145),
    @('ECs', 'Edn1', 'Ednrb', 'Inflammatory-Mac', 3, 1, 20.66988266666667, 62.009648, 0.9592026714402365, 0.9592026714402365, 3, 1, 2.195310666666666, 6.585932, 0.07524896398496668, 0.07524896398496668, 45.37681389688177, 408.391325071936, 0.07217900727749019, 0.07217900727749019),
    @('ECs', 'Edn1', 'Ednrb', 'MuSCs', 3, 1, 20.66988266666667, 62.009648, 0.9592026714402365, 0.9592026714402365, 3, 1, 3.329509666666667, 9.988529, 0.1141260582380437, 0.1141260582380437, 68.82057414753244, 619.3851673277919, 0.1094700199428755, 0.1094700199428755),
    @('ECs', 'Edn1', 'Ednrb', 'Resolving-Mac', 3, 1, 20.66988266666667, 62.009648, 0.9592026714402365, 0.9592026714402365, 3, 1, 2.438488333333333, 7.315465, 0.08358439812592726, 0.08358439812592725, 50.40326773403555, 453.62940960632, 0.08017437797311372, 0.08017437797311371),
    @('FAPs', 'Edn1', 'Ednrb', 'ECs', 2, 0.6666666666666666, 0.232822, 0.698466, 0.01080429376264442, 0.01080429376264442, 3, 1, 20.92735966666666, 62.782079, 0.717329969634113, 0.717329969634113, 4.872349732312666, 43.851147590814, 0.00775024371667576, 0.007750243716675762),
    @('FAPs', 'Edn1', 'Ednrb', 'FAPs', 2, 0.6666666666666666, 0.232822, 0.698466, 0.01080429376264442, 0.01080429376264442, 1, 0.3333333333333333, 0.283297, 0.849891, 0.009710610016949358, 0.009710610016949358, 0.06595777413399999, 0.593619967206, 0.0001049162832375984, 0.0001049162832375984),
    @('FAPs', 'Edn1', 'Ednrb', 'Inflammatory-Mac', 2, 0.6666666666666666, 0.232822, 0.698466, 0.01080429376264442, 0.01080429376264442, 3, 1, 2.195310666666666, 6.585932, 0.07524896398496668, 0.07524896398496668, 0.5111166200346666, 4.600049580312, 0.0008130119122282303, 0.0008130119122282304),
    @('FAPs', 'Edn1', 'Ednrb', 'MuSCs', 2, 0.6666666666666666, 0.232822, 0.698466, 0.01080429376264442, 0.01080429376264442, 3, 1, 3.329509666666667, 9.988529, 0.1141260582380437, 0.1141260582380437, 0.7751830996126666, 6.976647896514, 0.001233051459176489, 0.001233051459176489),
    @('FAPs', 'Edn1', 'Ednrb', 'Resolving-Mac', 2, 0.6666666666666666, 0.232822, 0.698466, 0.01080429376264442, 0.01080429376264442, 3, 1, 2.438488333333333, 7.315465, 0.08358439812592726, 0.08358439812592725, 0.5677337307433333, 5.10960357669, 0.0009030703913263441, 0.0009030703913263441),
    @('Inflammatory-Mac', 'Edn1', 'Ednrb', 'ECs', 2, 0.6666666666666666, 0.144358, 0.433074, 0.006699050085420723, 0.006699050085420723, 3, 1, 20.92735966666666, 62.782079, 0.717329969634113, 0.717329969634113, 3.021031786760666, 27.189286080846, 0.00480542939435225, 0.00480542939435225),
    @('Inflammatory-Mac', 'Edn1', 'Ednrb', 'FAPs', 2, 0.6666666666666666, 0.144358, 0.433074, 0.006699050085420723, 0.006699050085420723, 1, 0.3333333333333333, 0.283297, 0.849891, 0.009710610016949358, 0.009710610016949358, 0.04089618832599999, 0.3680656949339999, 0.00006505186286353193, 0.00006505186286353193),
    @('Inflammatory-Mac', 'Edn1', 'Ednrb', 'Inflammatory-Mac', 2, 0.6666666666666666, 0.144358, 0.433074, 0.006699050085420723, 0.006699050085420723, 3, 1, 2.195310666666666, 6.585932, 0.07524896398496668, 0.07524896398496668, 0.3169106572186666, 2.852195914968, 0.000504096578611312, 0.000504096578611312),
    @('Inflammatory-Mac', 'Edn1', 'Ednrb', 'MuSCs', 2, 0.6666666666666666, 0.144358, 0.433074, 0.006699050085420723, 0.006699050085420723, 3, 1, 3.329509666666667, 9.988529, 0.1141260582380437, 0.1141260582380437, 0.4806413564606666, 4.325772208146, 0.0007645361801882968, 0.0007645361801882968),
    @('Inflammatory-Mac', 'Edn1', 'Ednrb', 'Resolving-Mac', 2, 0.6666666666666666, 0.144358, 0.433074, 0.006699050085420723, 0.006699050085420723, 3, 1, 2.438488333333333, 7.315465, 0.08358439812592726, 0.08358439812592725, 0.3520152988233333, 3.16813768941, 0.0005599360694053327, 0.0005599360694053326),
    @('MuSCs', 'Edn1', 'Ednrb', 'ECs', 3, 1, 0.192848, 0.5785439999999999, 0.008949267867892432, 0.008949267867892432, 3, 1, 20.92735966666666, 62.782079, 0.717329969634113, 0.717329969634113, 4.035799456997332, 36.322195112976, 0.006419578047922822, 0.006419578047922822),
    @('MuSCs', 'Edn1', 'Ednrb', 'FAPs', 3, 1, 0.192848, 0.5785439999999999, 0.008949267867892432, 0.008949267867892432, 1, 0.3333333333333333, 0.283297, 0.849891, 0.009710610016949358, 0.009710610016949358, 0.05463325985599999, 0.4916993387039999, 0.00008690285020231928, 0.00008690285020231928),
    @('MuSCs', 'Edn1', 'Ednrb', 'Inflammatory-Mac', 3, 1, 0.192848, 0.5785439999999999, 0.008949267867892432, 0.008949267867892432, 3, 1, 2.195310666666666, 6.585932, 0.07524896398496668, 0.07524896398496668, 0.4233612714453333, 3.810251443008, 0.0006734231354828572, 0.0006734231354828572),
    @('MuSCs', 'Edn1', 'Ednrb', 'MuSCs', 3, 1, 0.192848, 0.5785439999999999, 0.008949267867892432, 0.008949267867892432, 3, 1, 3.329509666666667, 9.988529, 0.1141260582380437, 0.1141260582380437, 0.6420892801973332, 5.778803521775999, 0.001021344665878944, 0.001021344665878944),
    @('MuSCs', 'Edn1', 'Ednrb', 'Resolving-Mac', 3, 1, 0.192848, 0.5785439999999999, 0.008949267867892432, 0.008949267867892432, 3, 1, 2.438488333333333, 7.315465, 0.08358439812592726, 0.08358439812592725, 0.4702575981066667, 4.232318382959999, 0.0007480191684054893, 0.0007480191684054891),
    @('Resolving-Mac', 'Edn1', 'Ednrb', 'ECs', 3, 1, 0.3091146666666666, 0.9273439999999999, 0.0143447168438059, 0.0143447168438059, 3, 1, 20.92735966666666, 62.782079, 0.717329969634113, 0.717329969634113, 6.46895380757511, 58.22058426817599, 0.01028989529797723, 0.01028989529797723),
    @('Resolving-Mac', 'Edn1', 'Ednrb', 'FAPs', 3, 1, 0.3091146666666666, 0.9273439999999999, 0.0143447168438059, 0.0143447168438059, 1, 0.3333333333333333, 0.283297, 0.849891, 0.009710610016949358, 0.009710610016949358, 0.08757125772266665, 0.7881413195039999, 0.0001392959510737637, 0.0001392959510737637),
    @('Resolving-Mac', 'Edn1', 'Ednrb', 'Inflammatory-Mac', 3, 1, 0.3091146666666666, 0.9273439999999999, 0.0143447168438059, 0.0143447168438059, 3, 1, 2.195310666666666, 6.585932, 0.07524896398496668, 0.07524896398496668, 0.6786027249564444, 6.107424524608, 0.001079425081154095, 0.001079425081154095),
    @('Resolving-Mac', 'Edn1', 'Ednrb', 'MuSCs', 3, 1, 0.3091146666666666, 0.9273439999999999, 0.0143447168438059, 0.0143447168438059, 3, 1, 3.329509666666667, 9.988529, 0.1141260582380437, 0.1141260582380437, 1.029200270775111, 9.262802436975999, 0.001637105989924438, 0.001637105989924438),
    @('Resolving-Mac', 'Edn1', 'Ednrb', 'Resolving-Mac', 3, 1, 0.3091146666666666, 0.9273439999999999, 0.0143447168438059, 0.0143447168438059, 3, 1, 2.438488333333333, 7.315465, 0.08358439812592726, 0.08358439812592725, 0.7537725083288889, 6.783952574959999, 0.001198994523676367, 0.001198994523676367)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
